$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.364128232002258
$ws.Range("B1").Value = 2.849937677383423
$ws.Range("C1").Value = 3.85257887840271
$ws.Range("D1").Value = 3.518016815185547
$ws.Range("E1").Value = 0.8124356865882874
